$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2033898305084746
$ws.Range("C2").Value = 0.5694915254237288
$ws.Range("J2").Value = 0.02033898305084746
$ws.Range("P2").Value = 0.1389830508474576
$ws.Range("S2").Value = 0.06779661016949153
$ws.Range("C3").Value = 0.04597701149425287
$ws.Range("J3").Value = 0.01149425287356322
$ws.Range("O3").Value = 0.005747126436781609
$ws.Range("P3").Value = 0.7528735632183908
$ws.Range("S3").Value = 0.1839080459770115
$ws.Range("P4").Value = 0.7073170731707317
$ws.Range("S4").Value = 0.2926829268292683
$ws.Range("P5").Value = 0.3333333333333333
$ws.Range("S5").Value = 0.6666666666666666
$ws.Range("B6").Value = 0.05855855855855856
$ws.Range("E6").Value = 0.004504504504504504
$ws.Range("F6").Value = 0.07657657657657657
$ws.Range("J6").Value = 0.2297297297297297
$ws.Range("O6").Value = 0.03153153153153153
$ws.Range("Q6").Value = 0.1081081081081081
$ws.Range("R6").Value = 0.08108108108108109
$ws.Range("S6").Value = 0.4099099099099099
$ws.Range("B7").Value = 0.0989010989010989
$ws.Range("D7").Value = 0.02197802197802198
$ws.Range("E7").Value = 0.01098901098901099
$ws.Range("F7").Value = 0.05494505494505494
$ws.Range("J7").Value = 0.1043956043956044
$ws.Range("O7").Value = 0.005494505494505495
$ws.Range("Q7").Value = 0.1318681318681319
$ws.Range("R7").Value = 0.1043956043956044
$ws.Range("S7").Value = 0.467032967032967
$ws.Range("B8").Value = 0.08905380333951762
$ws.Range("D8").Value = 0.01484230055658627
$ws.Range("F8").Value = 0.04081632653061224
$ws.Range("J8").Value = 0.09090909090909091
$ws.Range("O8").Value = 0.005565862708719851
$ws.Range("Q8").Value = 0.1651205936920223
$ws.Range("R8").Value = 0.1205936920222635
$ws.Range("S8").Value = 0.4730983302411874
$ws.Range("B9").Value = 0.09405940594059406
$ws.Range("D9").Value = 0.009900990099009901
$ws.Range("F9").Value = 0.0594059405940594
$ws.Range("J9").Value = 0.1237623762376238
$ws.Range("O9").Value = 0.004950495049504951
$ws.Range("Q9").Value = 0.1435643564356436
$ws.Range("R9").Value = 0.0891089108910891
$ws.Range("S9").Value = 0.4752475247524752
$ws.Range("B10").Value = 0.09659090909090909
$ws.Range("D10").Value = 0.01917613636363636
$ws.Range("E10").Value = 0.001420454545454545
$ws.Range("F10").Value = 0.07315340909090909
$ws.Range("J10").Value = 0.1022727272727273
$ws.Range("O10").Value = 0.009943181818181818
$ws.Range("Q10").Value = 0.1938920454545454
$ws.Range("R10").Value = 0.09232954545454546
$ws.Range("S10").Value = 0.4112215909090909
$ws.Range("G11").Value = 0.1470588235294118
$ws.Range("J11").Value = 0.07352941176470588
$ws.Range("K11").Value = 0.1764705882352941
$ws.Range("L11").Value = 0.5919117647058824
$ws.Range("S11").Value = 0.01102941176470588
$ws.Range("G12").Value = 0.7218934911242604
$ws.Range("J12").Value = 0.1952662721893491
$ws.Range("L12").Value = 0.03550295857988166
$ws.Range("S12").Value = 0.04733727810650887
$ws.Range("F13").Value = 0.02564102564102564
$ws.Range("G13").Value = 0.6153846153846154
$ws.Range("J13").Value = 0.2307692307692308
$ws.Range("S13").Value = 0.1282051282051282
$ws.Range("F15").Value = 0.01895734597156398
$ws.Range("H15").Value = 0.1943127962085308
$ws.Range("I15").Value = 0.06161137440758294
$ws.Range("J15").Value = 0.4028436018957346
$ws.Range("K15").Value = 0.06161137440758294
$ws.Range("M15").Value = 0.004739336492890996
$ws.Range("N15").Value = 0.004739336492890996
$ws.Range("O15").Value = 0.05213270142180094
$ws.Range("S15").Value = 0.1990521327014218
$ws.Range("F16").Value = 0.01538461538461539
$ws.Range("H16").Value = 0.1948717948717949
$ws.Range("I16").Value = 0.04615384615384616
$ws.Range("J16").Value = 0.4307692307692308
$ws.Range("K16").Value = 0.1435897435897436
$ws.Range("M16").Value = 0.01538461538461539
$ws.Range("N16").Value = 0.005128205128205128
$ws.Range("O16").Value = 0.04615384615384616
$ws.Range("S16").Value = 0.1025641025641026
$ws.Range("F17").Value = 0.009029345372460496
$ws.Range("H17").Value = 0.2031602708803612
$ws.Range("I17").Value = 0.1038374717832957
$ws.Range("J17").Value = 0.4469525959367946
$ws.Range("K17").Value = 0.08126410835214447
$ws.Range("M17").Value = 0.004514672686230248
$ws.Range("N17").Value = 0.002257336343115124
$ws.Range("O17").Value = 0.04514672686230248
$ws.Range("S17").Value = 0.1038374717832957
$ws.Range("F18").Value = 0.008032128514056224
$ws.Range("H18").Value = 0.2248995983935743
$ws.Range("I18").Value = 0.0963855421686747
$ws.Range("J18").Value = 0.4136546184738956
$ws.Range("K18").Value = 0.06827309236947791
$ws.Range("M18").Value = 0.01606425702811245
$ws.Range("O18").Value = 0.08433734939759036
$ws.Range("S18").Value = 0.08835341365461848
$ws.Range("F19").Value = 0.01588397790055249
$ws.Range("H19").Value = 0.2147790055248619
$ws.Range("I19").Value = 0.07803867403314917
$ws.Range("J19").Value = 0.417817679558011
$ws.Range("K19").Value = 0.08701657458563536
$ws.Range("M19").Value = 0.02071823204419889
$ws.Range("N19").Value = 0.0006906077348066298
$ws.Range("O19").Value = 0.06146408839779006
$ws.Range("S19").Value = 0.1035911602209945
